# Update Leve profit-tracking values across all item sheets.
# These numbers (currentAveragePrice / NQ / HQ, LevePrice NQ/HQ, LeveProfit
# NQ/HQ) come from a scheduled market-data refresh; this mirrors that
# refreshed snapshot back into the workbook.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 850
$ws.Range("I2").Value = 850
$ws.Range("K2").Value = 850
$ws.Range("M2").Value = -737
$ws.Range("H9").Value = 108.6
$ws.Range("I9").Value = 108.6
$ws.Range("K9").Value = 108.6
$ws.Range("M9").Value = 60.40000000000001
$ws.Range("H51").Value = 6046.5
$ws.Range("I51").Value = 7595.6665
$ws.Range("J51").Value = 1399
$ws.Range("K51").Value = 7595.6665
$ws.Range("L51").Value = 1399
$ws.Range("M51").Value = -7111.6665
$ws.Range("N51").Value = -2367
$ws.Range("H135").Value = 17745.834
$ws.Range("J135").Value = 34666.332
$ws.Range("L135").Value = 311996.988
$ws.Range("N135").Value = -317066.988
$ws.Range("H137").Value = 4002.862
$ws.Range("I137").Value = 1913.1765
$ws.Range("K137").Value = 5739.529500000001
$ws.Range("M137").Value = -3189.529500000001

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1057
$ws.Range("I2").Value = 1057
$ws.Range("K2").Value = 1057
$ws.Range("M2").Value = -944
$ws.Range("H32").Value = 17244504
$ws.Range("I32").Value = 20001892
$ws.Range("K32").Value = 20001892
$ws.Range("M32").Value = -20001605
$ws.Range("H61").Value = 19278590
$ws.Range("I61").Value = 27785368
$ws.Range("K61").Value = 27785368
$ws.Range("M61").Value = -27785156
$ws.Range("H74").Value = 12514121
$ws.Range("I74").Value = 41670160
$ws.Range("J74").Value = 18677.143
$ws.Range("K74").Value = 41670160
$ws.Range("L74").Value = 18677.143
$ws.Range("M74").Value = -41669286
$ws.Range("N74").Value = -20425.143
$ws.Range("H77").Value = 12514121
$ws.Range("I77").Value = 41670160
$ws.Range("J77").Value = 18677.143
$ws.Range("K77").Value = 208350800
$ws.Range("L77").Value = 93385.715
$ws.Range("M77").Value = -208346432
$ws.Range("N77").Value = -102121.715
$ws.Range("H116").Value = 1057
$ws.Range("I116").Value = 1057
$ws.Range("K116").Value = 1057
$ws.Range("M116").Value = 1237
$ws.Range("H136").Value = 19278590
$ws.Range("I136").Value = 27785368
$ws.Range("K136").Value = 83356104
$ws.Range("M136").Value = -83353554

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1057
$ws.Range("I3").Value = 1057
$ws.Range("K3").Value = 1057
$ws.Range("M3").Value = -943
$ws.Range("H86").Value = 2400.8
$ws.Range("I86").Value = 2256.4443
$ws.Range("K86").Value = 2256.4443
$ws.Range("M86").Value = -1133.4443
$ws.Range("H89").Value = 2400.8
$ws.Range("I89").Value = 2256.4443
$ws.Range("K89").Value = 11282.2215
$ws.Range("M89").Value = -5666.2215
$ws.Range("H134").Value = 92370.086
$ws.Range("J134").Value = 267507
$ws.Range("L134").Value = 802521
$ws.Range("N134").Value = -807591

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 2802.923
$ws.Range("I99").Value = 3356
$ws.Range("J99").Value = 2557.111
$ws.Range("K99").Value = 3356
$ws.Range("L99").Value = 2557.111
$ws.Range("M99").Value = -1858
$ws.Range("N99").Value = -5553.111
$ws.Range("H126").Value = 2802.923
$ws.Range("I126").Value = 3356
$ws.Range("J126").Value = 2557.111
$ws.Range("K126").Value = 10068
$ws.Range("L126").Value = 7671.333
$ws.Range("M126").Value = -7598
$ws.Range("N126").Value = -12611.333
$ws.Range("H134").Value = 595847.4
$ws.Range("I134").Value = 837197.8
$ws.Range("K134").Value = 2511593.4
$ws.Range("M134").Value = -2509058.4

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 453562.44
$ws.Range("I12").Value = 1912.8182
$ws.Range("K12").Value = 5738.4546
$ws.Range("M12").Value = -5565.4546
$ws.Range("H119").Value = 3535.75
$ws.Range("I119").Value = 1242.9
$ws.Range("K119").Value = 3728.7
$ws.Range("M119").Value = 1109.3

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 621.8
$ws.Range("I14").Value = 402.66666
$ws.Range("K14").Value = 402.66666
$ws.Range("M14").Value = -234.66666
$ws.Range("H111").Value = 57424.5
$ws.Range("J111").Value = 57424.5
$ws.Range("L111").Value = 57424.5
$ws.Range("N111").Value = -63558.5
$ws.Range("H122").Value = 2719.4
$ws.Range("I122").Value = 2899.25
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 8697.75
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -6247.75
$ws.Range("N122").Value = -10900
$ws.Range("H124").Value = 115000
$ws.Range("J124").Value = 115000
$ws.Range("L124").Value = 115000
$ws.Range("N124").Value = -124820
$ws.Range("H126").Value = 3671
$ws.Range("I126").Value = 3741.5715
$ws.Range("K126").Value = 11224.7145
$ws.Range("M126").Value = -8754.7145
$ws.Range("H132").Value = 52647156
$ws.Range("I132").Value = 71434350
$ws.Range("J132").Value = 42999.8
$ws.Range("K132").Value = 214303050
$ws.Range("L132").Value = 128999.4
$ws.Range("M132").Value = -214300520
$ws.Range("N132").Value = -134059.4

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 943.4375
$ws.Range("I61").Value = 947.7273
$ws.Range("J61").Value = 934
$ws.Range("K61").Value = 947.7273
$ws.Range("L61").Value = 934
$ws.Range("M61").Value = -745.7273
$ws.Range("N61").Value = -1338
$ws.Range("H81").Value = 58330
$ws.Range("I81").Value = 50000
$ws.Range("J81").Value = 62495
$ws.Range("K81").Value = 50000
$ws.Range("L81").Value = 62495
$ws.Range("M81").Value = -49002
$ws.Range("N81").Value = -64491
$ws.Range("H84").Value = 58330
$ws.Range("I84").Value = 50000
$ws.Range("J84").Value = 62495
$ws.Range("K84").Value = 150000
$ws.Range("L84").Value = 187485
$ws.Range("M84").Value = -145008
$ws.Range("N84").Value = -197469
$ws.Range("H92").Value = 134000
$ws.Range("J92").Value = 134000
$ws.Range("L92").Value = 134000
$ws.Range("N92").Value = -138992
$ws.Range("H113").Value = 943.4375
$ws.Range("I113").Value = 947.7273
$ws.Range("J113").Value = 934
$ws.Range("K113").Value = 947.7273
$ws.Range("L113").Value = 934
$ws.Range("M113").Value = 1222.2727
$ws.Range("N113").Value = -5274
$ws.Range("H132").Value = 53716.168
$ws.Range("I132").Value = 6508.7144
$ws.Range("K132").Value = 19526.1432
$ws.Range("M132").Value = -16996.1432

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H125").Value = 74049.5
$ws.Range("J125").Value = 74049.5
$ws.Range("L125").Value = 74049.5
$ws.Range("N125").Value = -83889.5

